# Adds a new "2022-Q3" quarter sheet (right after the "总计" summary sheet)
# and inserts its corresponding summary row into "总计".
#
# Net effect on sheet order:
#   总计, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3, 2021-Q1, 2020-Q4
#   -> 总计, 2022-Q3, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3, 2021-Q1, 2020-Q4

$wb = $excel.ActiveWorkbook

# Remember which sheet was originally the active/selected tab so we can
# restore that UI state after inserting the new sheet (Add() makes the new
# sheet active, which we don't want to leave behind).
$originalActiveSheetName = $wb.Worksheets.Item($wb.Worksheets.Count).Name

function Set-TextValue {
    param($rng, $val)
    # Force the cell to be stored as text (not auto-converted to a number)
    # without leaving a lingering custom number-format style behind.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# ---------------------------------------------------------------------------
# 1) Insert the new summary row into "总计" (sheet 1), right under the header.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.21

# A2 needs the same style as the other index cells (A3 still has it after
# the insert) plus its own sequential index value.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("A2").Value = 0

# Renumber the index column for all the rows that shifted down by one.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6

# ---------------------------------------------------------------------------
# 2) Add the new "2022-Q3" sheet right after "总计" (i.e. before the sheet
#    that is currently in position 2, "2022-Q2").
# ---------------------------------------------------------------------------
$nextSheet = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($nextSheet, $null)
$q3.Name = "2022-Q3"

# Match the outline settings used by the other quarter sheets
# (<outlinePr summaryBelow="1" summaryRight="1"/>).
$q3.Outline.SummaryRow = 1
$q3.Outline.SummaryColumn = 1

# Grab formatting references from the (now-shifted) "2022-Q2" sheet so the
# new sheet's header row / index column match the existing look (style "2":
# bold, bordered, centered).
$q2 = $wb.Worksheets.Item(3)

$q2.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q2.Range("A2:A5").Copy()
$q3.Range("A2:A5").PasteSpecial(-4122)
$q3.Range("A2").Value = 0
$q3.Range("A3").Value = 1
$q3.Range("A4").Value = 2
$q3.Range("A5").Value = 3

Set-TextValue $q3.Range("B2") "005994"
Set-TextValue $q3.Range("C2") "国投瑞银中证500指数量化增强A"
Set-TextValue $q3.Range("D2") "13.36"
Set-TextValue $q3.Range("E2") "88.67"
Set-TextValue $q3.Range("F2") "1.13"
Set-TextValue $q3.Range("G2") "0.1510"
$q3.Range("H2").Value = 9

Set-TextValue $q3.Range("B3") "007089"
Set-TextValue $q3.Range("C3") "国投瑞银中证500指数量化增强C"
Set-TextValue $q3.Range("D3") "4.45"
Set-TextValue $q3.Range("E3") "88.67"
Set-TextValue $q3.Range("F3") "1.13"
Set-TextValue $q3.Range("G3") "0.0503"
$q3.Range("H3").Value = 9

Set-TextValue $q3.Range("B4") "015889"
Set-TextValue $q3.Range("C4") "富国中证500基本面精选股票A"
Set-TextValue $q3.Range("D4") "0.97"
Set-TextValue $q3.Range("E4") "40.38"
Set-TextValue $q3.Range("F4") "0.62"
Set-TextValue $q3.Range("G4") "0.0060"
$q3.Range("H4").Value = 8

Set-TextValue $q3.Range("B5") "015890"
Set-TextValue $q3.Range("C5") "富国中证500基本面精选股票C"
Set-TextValue $q3.Range("D5") "0.06"
Set-TextValue $q3.Range("E5") "40.38"
Set-TextValue $q3.Range("F5") "0.62"
Set-TextValue $q3.Range("G5") "0.0004"
$q3.Range("H5").Value = 8

# ---------------------------------------------------------------------------
# 3) Restore the originally-active sheet tab (Add() switched focus to the
#    new sheet).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item($originalActiveSheetName).Activate()
